$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Implementation Details" paragraph: merge the run split around
#    "PostgreSQL" (removes the spell-check proofErr markers as a
#    natural side effect of replacing text that spans them).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "database using PostgreSQL as the model",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "database using PostgreSQL as the model", 2)

# ------------------------------------------------------------------
# 2. "Security" paragraph: merge the run split around "Csrf" (removes
#    the spell-check proofErr markers around it).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "any potential SQL injections.   Csrf tokens",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "any potential SQL injections.   Csrf tokens", 2)

# ------------------------------------------------------------------
# 3. "Performance" paragraph: merge the run split around "NoSQL"
#    (removes the spell-check proofErr markers around it).  The
#    match purposely starts at "NoSQL" (rather than earlier in the
#    sentence) so the leading "<w:tab/>" of the paragraph's first run
#    is left completely untouched.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "NoSQL to prevent",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NoSQL to prevent", 2)

# ------------------------------------------------------------------
# 4. "Challenges" paragraph: "qualify" -> "qualification" and drop
#    the trailing "  There " left dangling at the end of the section.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "qualify since the law for every country.  There ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "qualification since the law for every country.", 2)

# ------------------------------------------------------------------
# 5. Move the "_GoBack" bookmark (tracks the location of the most
#    recent edit) from the "security" paragraph to the end of the
#    sentence we just edited in the "Challenges" paragraph.
#
#    A temporary padding string is appended at the very end of the
#    document first: placing a zero-length bookmark on the last
#    couple of character positions of the document is unreliable, so
#    we push our target position away from the document's tail,
#    create the bookmark, then remove the padding again.
# ------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("ZZ_TEMP_PADDING_ZZ")

$targetRange = $d.Content
$targetRange.Find.Execute("law for every country.")
$targetRange.Collapse(0)

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$d.Bookmarks.Add("_GoBack", $targetRange)

$padRange = $d.Range($d.Content.End - 19, $d.Content.End - 1)
$padRange.Delete()
